$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.313.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.865.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.46%  "

$ws.Range("E4").Value = "  +0.56%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.49%  "

$ws.Range("E6").Value = "  +0.90%  "

$ws.Range("E7").Value = "  +0.53%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.68"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.96%  "

$ws.Range("E9").Value = "  +0.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0697"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0990"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.80%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.135.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.33%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.866.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.48%  "

$ws.Range("E15").Value = "  +1.18%  "

$ws.Range("E16").Value = "  +1.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.276.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0799"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "241.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.29%  "

$ws.Range("E21").Value = "  +0.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.49%  "

$ws.Range("E23").Value = "  +0.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "169.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +25.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.19"
$ws.Range("D27").Style = "Normal"

$ws.Range("E28").Value = "  +1.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0564"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.06%  "

$ws.Range("E32").Value = "  +2.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +28.24%  "

$ws.Range("E34").Value = "  +2.37%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.77%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.818"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +17.38%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.88%  "

$ws.Range("E38").Value = "  +4.58%  "

$ws.Range("E39").Value = "  +4.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "90.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.43%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.346.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.07%  "

$ws.Range("E42").Value = "  +1.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0601"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +15.03%  "

$ws.Range("E45").Value = "  +0.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +48.50%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.31%  "

$ws.Range("E48").Value = "  -0.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.051.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0686"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.22%  "

$ws.Range("E51").Value = "  +0.59%  "

